{"js": "// Update the King contract's deployed address and wrap the new value in a\n// bookmark (mirrors the \"Insert Bookmark\" step Word performs when the\n// address run is selected and tagged, e.g. via a DDE/quick-link paste).\n\nconst OLD_ADDRESS = \"0x0d6f470dada8f8390046c01015f8924e5055ac54\";\nconst NEW_ADDRESS = \"0x65840e9c5dbbf36a7aed1cce4893b2f1218bcd6a\";\nconst BOOKMARK_NAME = \"__DdeLink__104_2055276817\";\n\nconst body = context.document.body;\n\n// Locate the run containing the old contract address and swap its text.\nconst matches = body.search(OLD_ADDRESS, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Contract address run not found: \" + OLD_ADDRESS);\n}\n\nconst addressRange = matches.items[0];\naddressRange.insertText(NEW_ADDRESS, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate the (now updated) address text and wrap it in a bookmark so the\n// run is bracketed by bookmarkStart/bookmarkEnd exactly as in the source.\nconst newMatches = body.search(NEW_ADDRESS, { matchCase: true });\nnewMatches.load(\"items\");\nawait context.sync();\n\nif (newMatches.items.length === 0) {\n  throw new Error(\"Updated contract address run not found: \" + NEW_ADDRESS);\n}\n\nnewMatches.items[0].insertBookmark(BOOKMARK_NAME);\nawait context.sync();\n", "ps1": "# Update the King contract's deployed address and wrap the new value in a\n# bookmark (mirrors the \"Insert Bookmark\" step Word performs when the\n# address run is selected and tagged, e.g. via a DDE/quick-link paste).\n\n$d = $word.ActiveDocument\n\n$oldAddress = \"0x0d6f470dada8f8390046c01015f8924e5055ac54\"\n$newAddress = \"0x65840e9c5dbbf36a7aed1cce4893b2f1218bcd6a\"\n$bookmarkName = \"__DdeLink__104_2055276817\"\n\n# Locate the run containing the old contract address.\n$found = $d.Content\n$found.Find.ClearFormatting()\n$found.Find.Text = $oldAddress\n$found.Find.Forward = $true\n$ok = $found.Find.Execute()\nif (-not $ok) {\n    throw \"Could not find contract address run\"\n}\n$runStart = $found.Start\n$runEnd = $found.End\n\n# Replace the address text while preserving the run's character formatting\n# (italic Times New Roman). Re-pointing Range.Text at a span whose Start is\n# exactly the run's own start resets that run's rPr, so keep at least the\n# shared leading characters (\"0x\") anchored and only rewrite the remainder.\n$prefixLen = 0\n$maxLen = [Math]::Min($oldAddress.Length, $newAddress.Length)\nwhile ($prefixLen -lt $maxLen -and $oldAddress[$prefixLen] -eq $newAddress[$prefixLen]) {\n    $prefixLen++\n}\nif ($prefixLen -eq 0) { $prefixLen = 1 }\n\n$editRange = $d.Range($runStart + $prefixLen, $runEnd)\n$editRange.Text = $newAddress.Substring($prefixLen)\n\n# Re-select the (now updated) address text and wrap it in a bookmark so the\n# run is bracketed by bookmarkStart/bookmarkEnd exactly as in the source.\n$bmRange = $d.Range($runStart, $runStart + $newAddress.Length)\n$d.Bookmarks.Add($bookmarkName, $bmRange)\n"}
